# "Generate Report for Archive"
#
# The localization-status report is regenerated: the handoff/translation
# status moves from "Ready for handoff" to "In Translation" on every sheet
# that tracks it (the Overview roll-up columns for zh-cn/de-de, plus each
# language sheet's own Status column). Excel's column autosize then
# narrows those status columns to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Row 2 holds the (only) data record on each sheet.
# Overview: column E = zh-cn status, column F = de-de status.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# zh-cn / de-de sheets: column C = Status.
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# The status columns auto-narrow to fit the shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
